$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Metadata sheet: bump the "Date" value (B8) to the new timestamp.
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# ------------------------------------------------------------------
# 2. Elements sheet: the two "Mapping" columns (AK / AL) swap places
#    - header text swaps
#    - every data row's content swaps between AK and AL
#    - column widths swap accordingly
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

function Swap-CellValues($sheet, $cellA, $cellB) {
    $a = $sheet.Range($cellA).Value2
    $b = $sheet.Range($cellB).Value2
    $sheet.Range($cellA).Value = $b
    $sheet.Range($cellB).Value = $a
}

# Header row
Swap-CellValues $ws "AK1" "AL1"

# Data rows (2 and 4 have no Mapping content, only 3, 5 and 6 do)
Swap-CellValues $ws "AK3" "AL3"
Swap-CellValues $ws "AK5" "AL5"
Swap-CellValues $ws "AK6" "AL6"

# Column widths swap too (AK was narrow/84.4->24.98, AL was 84.4, now reversed)
$ws.Columns.Item(37).ColumnWidth = 83.5
$ws.Columns.Item(38).ColumnWidth = 24.166666666666668
